# Apply odds updates to the "Jogos da Semana" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 3 ---
$ws.Range("K3").Value  = 1.91
$ws.Range("M3").Value  = 1.13
$ws.Range("N3").Value  = 6
$ws.Range("U3").Value  = 4.2
$ws.Range("V3").Value  = 1.22
$ws.Range("AF3").Value = 17
$ws.Range("AG3").Value = 21
$ws.Range("AI3").Value = 6
$ws.Range("AK3").Value = 21
$ws.Range("AN3").Value = 9

# --- Row 4 ---
$ws.Range("G4").Value  = 1.75
$ws.Range("H4").Value  = 3.5
$ws.Range("O4").Value  = 1.4
$ws.Range("P4").Value  = 2.75
$ws.Range("AD4").Value = 7.5
$ws.Range("AG4").Value = 17

# --- Row 5 ---
$ws.Range("J5").Value  = 4.55
$ws.Range("L5").Value  = 2.4
$ws.Range("O5").Value  = 1.23
$ws.Range("P5").Value  = 3.35
$ws.Range("S5").Value  = 1.7
$ws.Range("T5").Value  = 1.91
$ws.Range("W5").Value  = 2.6
$ws.Range("X5").Value  = 1.38
$ws.Range("AB5").Value = 2.07
$ws.Range("AC5").Value = 13.5
$ws.Range("AD5").Value = 27
$ws.Range("AG5").Value = 40
$ws.Range("AH5").Value = 37
$ws.Range("AN5").Value = 8
$ws.Range("AO5").Value = 9.5
$ws.Range("AQ5").Value = 16.5
$ws.Range("AR5").Value = 13.5
$ws.Range("AS5").Value = 21
